# Update countries & provincias Spain
# Applies the 29-Abril-2020 18:22 data refresh to the "Pais" sheet:
#  - Poland's case counts overtook Qatar's, so the two countries swap
#    ranking rows (row 34 / row 35); Poland gets fresh numbers while
#    Qatar's figures are carried over unchanged.
#  - Several other country rows receive updated totals.
#  - The "last updated" timestamp footer moves from 17:52 to 18:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 1040233 4468 143302 837112 19110 553 59819

# Italia
Set-Row 6 "Italia" 203591 2086 71252 104657 1795 323 27682

# Chile
Set-Row 31 "Chile" 14885 520 8057 6612 377 9 216

# Row 34 used to be Catar; Polonia's updated numbers now outrank it
Set-Row 34 "Polonia" 12640 422 3025 8991 160 28 624

# Row 35 used to be Polonia; Catar drops to this row, values unchanged
Set-Row 35 "Catar" 12564 643 1243 11311 72 0 10

# Chequia
Set-Row 45 "Chequia" 7563 59 3096 4240 71 0 227

# Egipto
Set-Row 52 "Egipto" 5268 226 1335 3553 0 21 380

# Luxemburgo
Set-Row 59 "Luxemburgo" 3769 28 3134 546 21 0 89

# Isla de Man
Set-Row 127 "Isla de Man" 313 4 252 40 22 0 21

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 18:22"
